$d = $word.ActiveDocument
$apos = [char]8217

# ------------------------------------------------------------------
# 1) Collapse the "SQLite" / "onUpgrade" spell-check run-splits back
#    into the surrounding run (drops the now-stray w:proofErr marks
#    Word inserted around those two words) while leaving the
#    following "Thirdly, ..." sentence as its own separate run.
# ------------------------------------------------------------------
$target = "There were two reasons this happened. First, we were not at all experienced with Android SQLite so we didn" + $apos + "t implement the onUpgrade() method properly, and didn" + $apos + "t know that it wasn" + $apos + "t implemented properly because we didn" + $apos + "t realize the consequences of leaving the database version unchanged after changing the implementation of the database tables. Second, we failed to test the application itself frequently enough to detect the problem. "

$rng = $d.Content
$rng.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, $target, 2)

# The replace above re-normalizes the whole paragraph's runs, which
# also swallows the following "Thirdly, ..." run into the same run.
# Re-split the paragraph right before "Thirdly" so it again becomes
# its own run, matching the target structure.
$para = $d.Paragraphs(3).Range
$paraText = $para.Text
$idx = $paraText.IndexOf("Thirdly")
if ($idx -ge 0) {
    $splitPoint = $para.Start + $idx
    $splitRng = $d.Range($splitPoint, $splitPoint)
    $splitRng.InsertParagraphAfter()
    # Remove the paragraph mark that was just inserted so the two
    # pieces of text stay in the same paragraph, just as two runs.
    $markRng = $d.Range($splitPoint, $splitPoint + 1)
    $markRng.Delete()
}

# ------------------------------------------------------------------
# 2) Append a new, wholly empty paragraph after the final paragraph
#    of body text (right before the sectPr).
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$tailText = "since iteration 2. "
$tailRng = $d.Content
$tailRng.Find.Execute($tailText, $true, $false, $false, $false, $false, $true, 1, $false, ($tailText + "^p"), 2)

Write-Output "done"
